# Mini corrections: BS no intermediate result; SO delivery interval
#
# The row for Solothurn (row 12) previously showed the short
# "every 5 minutes" interval texts. Correct them to the longer texts
# that already exist elsewhere in the workbook ("... sofern neue
# Gemeindeergebnisse vorhanden sind" / "... si de nouveaux resultats
# communales sont disponibles"), matching row 23's formatting
# (border/alignment style and taller wrapped row height).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the cell formatting (style/border) from the row that already uses
# the long-form text (row 23) onto row 12's C/D cells.
$ws.Range("C23:D23").Copy()
$ws.Range("C12:D12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the actual cell text.
$ws.Range("C12").Value = "Lieferintervall: alle 5 Minuten, sofern neue Gemeindeergebnisse vorhanden sind"
$ws.Range("D12").Value = "Intervalle de transfert: toutes les 5 minutes, si de nouveaux résultats communales sont disponibles"

# The row now wraps onto two lines, same as row 23.
$ws.Rows(12).RowHeight = 29.25

# Reflect the cursor/selection position saved with the workbook.
[void]$ws.Range("D12").Select()
